# The previous commit accidentally left Word's change-tracking turned on
# while comment authors were being edited, so each edited comment picked
# up a disambiguating " [n]" suffix on its author name (Word does this
# when it thinks multiple distinct "Peter C. Chapin" identities exist).
# Undo that: strip the " [<number>]" suffix from every comment whose
# author is "Peter C. Chapin [<n>]", restoring the plain "Peter C. Chapin".

$d = $word.ActiveDocument

foreach ($c in $d.Comments) {
    if ($c.Author -match "^Peter C\. Chapin \[\d+\]$") {
        $c.Author = "Peter C. Chapin"
    }
}
